$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 updates
$ws.Range("C5").Value = "Didn't get all warnings fixed and didn't improve architecture but got array of lights working"
$ws.Range("E5").Value = 6
$ws.Range("G5").Value = "Supervisor advised going through code line by line to find where warnings coming from.  Could be something to do with something not being set up correctly for the first model."

# Row 6 updates
$ws.Range("B6").Value = "Get all warnings fixed. Have member variable for device context in scene class as Get calls have an overhead. Start getting assets together"
$ws.Rows.Item(6).RowHeight = 75

# Row 7 updates
$ws.Range("B7").Value = "Have a scene setup with house and lamps. Extra: implement variance mapping"
$ws.Rows.Item(7).RowHeight = 45

# Update the selection to E5
$ws.Range("E5").Select()
